$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New alarm data rows to append (rows 12-17), matching the source export format.
# Columns A, B, C, E, F, G are all text in this sheet (even the numeric-looking
# ones); only D (Alarm number) is a real number.
$newRows = @(
    @{ Row = 12; A = $null;          B = $null;  C = $null;     D = 5; E = "Blood Pressure";    F = "Low";       G = "2024-10-30T16:39:38.332" },
    @{ Row = 13; A = "2024-11-04";   B = "111";  C = "AA111";   D = 6; E = "Blood Pressure";    F = "Very Low";  G = "2024-11-04T06:51:05.800" },
    @{ Row = 14; A = "2024-11-04";   B = "111";  C = "AA111";   D = 7; E = "Blood Pressure";    F = "Very High"; G = "2024-11-04T06:52:41.859" },
    @{ Row = 15; A = "2024-11-04";   B = "111";  C = "AA111";   D = 7; E = "Heart Rate";        F = "Static";    G = "2024-11-04T06:52:51.788" },
    @{ Row = 16; A = "2024-11-04";   B = "112";  C = "AA112";   D = 8; E = "Oxygen Saturation"; F = "High";      G = "2024-11-04T06:57:17.776" },
    @{ Row = 17; A = "2024-11-04";   B = "113";  C = "Aa1123";  D = 9; E = "Blood Pressure";    F = "Static";    G = "2024-11-04T07:01:58.019" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # A/B/C look like dates or plain numbers ("2024-11-04", "111", ...), so
    # without a text number format Excel would silently coerce them into a
    # date serial number / numeric value on assignment. Force text first.
    if ($null -ne $r.A) {
        $cell = $ws.Range("A$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r.A
    }
    if ($null -ne $r.B) {
        $cell = $ws.Range("B$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r.B
    }
    if ($null -ne $r.C) {
        $cell = $ws.Range("C$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $r.C
    }

    # D (Alarm number) is a genuine number in the source data.
    $ws.Range("D$rowNum").Value = $r.D

    # E/F hold plain vital-sign/value words (never numeric-looking), so they
    # can be assigned directly and keep the sheet's default formatting.
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F

    # G is a timestamp that includes a literal "T" separator, which Excel's
    # date/time parser does not recognize - but force text anyway to match
    # the source export (which stores every timestamp as text).
    $cell = $ws.Range("G$rowNum")
    $cell.NumberFormat = "@"
    $cell.Value = $r.G
}
